$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 394.5
$ws.Range("I2").Value = 116.666664
$ws.Range("J2").Value = 513.5714
$ws.Range("K2").Value = 116.666664
$ws.Range("L2").Value = 513.5714
$ws.Range("M2").Value = -3.666663999999997
$ws.Range("N2").Value = -739.5714

# Row 40
$ws.Range("H40").Value = 3300
$ws.Range("I40").Value = 4266.6665
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 4266.6665
$ws.Range("L40").Value = 2333.3333
$ws.Range("M40").Value = -4091.6665
$ws.Range("N40").Value = -2683.3333

# Row 121
$ws.Range("H121").Value = 1389.8334
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()

# Row 129
$ws.Range("H129").Value = 840.9756
$ws.Range("J129").Value = 1038.8334
$ws.Range("L129").Value = 3116.5002
$ws.Range("N129").Value = -13116.5002

# Row 138
$ws.Range("H138").Value = 1396.01
$ws.Range("J138").Value = 1531.8861
$ws.Range("L138").Value = 4595.6583
$ws.Range("N138").Value = -14875.6583

# Row 141
$ws.Range("H141").Value = 532.6923
$ws.Range("I141").Value = 551.75
$ws.Range("K141").Value = 1655.25
$ws.Range("M141").Value = 3524.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2454.7856
$ws.Range("I20").Value = 2538.7778
$ws.Range("J20").Value = 2303.6
$ws.Range("K20").Value = 2538.7778
$ws.Range("L20").Value = 2303.6
$ws.Range("M20").Value = -2291.7778
$ws.Range("N20").Value = -2797.6

# Row 22
$ws.Range("H22").Value = 1750
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1327
$ws.Range("N22").Value = -2346

# Row 81
$ws.Range("H81").Value = 10089.5
$ws.Range("J81").Value = 10089.5
$ws.Range("L81").Value = 10089.5
$ws.Range("N81").Value = -12211.5

# Row 84
$ws.Range("H84").Value = 10089.5
$ws.Range("J84").Value = 10089.5
$ws.Range("L84").Value = 30268.5
$ws.Range("N84").Value = -40876.5

# Row 86
$ws.Range("H86").Value = 3910.9473
$ws.Range("I86").Value = 4486.2856
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 4486.2856
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -3363.2856
$ws.Range("N86").Value = -4546

# Row 89
$ws.Range("H89").Value = 3910.9473
$ws.Range("I89").Value = 4486.2856
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 22431.428
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -16815.428
$ws.Range("N89").Value = -22732

# Row 94
$ws.Range("H94").Value = 16667466
$ws.Range("I94").Value = 16667466
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 16667466
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -16667015
$ws.Range("N94").ClearContents()

# Row 105
$ws.Range("H105").Value = 62501590
$ws.Range("I105").Value = 66668332
$ws.Range("J105").Value = 500
$ws.Range("K105").Value = 66668332
$ws.Range("L105").Value = 500
$ws.Range("M105").Value = -66666585
$ws.Range("N105").Value = -3994

# Row 107
$ws.Range("H107").Value = 1400.7142
$ws.Range("J107").Value = 966.6667
$ws.Range("L107").Value = 966.6667
$ws.Range("N107").Value = -4806.6667

# Row 137
$ws.Range("H137").Value = 38559
$ws.Range("J137").Value = 38559
$ws.Range("L137").Value = 38559
$ws.Range("N137").Value = -48759

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 54264.617
$ws.Range("I22").Value = 488.33334
$ws.Range("K22").Value = 488.33334
$ws.Range("M22").Value = -138.33334

# Row 31
$ws.Range("H31").Value = 1884
$ws.Range("I31").Value = 1919.36
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 1919.36
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -1624.36
$ws.Range("N31").Value = -1590

# Row 34
$ws.Range("H34").Value = 1884
$ws.Range("I34").Value = 1919.36
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1919.36
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -1717.36
$ws.Range("N34").Value = -1404

# Row 58
$ws.Range("H58").Value = 634.44446
$ws.Range("J58").Value = 691.6111
$ws.Range("L58").Value = 691.6111
$ws.Range("N58").Value = -1097.6111

# Row 99
$ws.Range("H99").Value = 1921.4667
$ws.Range("I99").Value = 1955.5385
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1955.5385
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = -457.5385000000001
$ws.Range("N99").Value = -4696

# Row 126
$ws.Range("H126").Value = 1921.4667
$ws.Range("I126").Value = 1955.5385
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 5866.6155
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -3396.6155
$ws.Range("N126").Value = -10040

# Row 132
$ws.Range("H132").Value = 3466.7778
$ws.Range("I132").Value = 3108
$ws.Range("K132").Value = 9324
$ws.Range("M132").Value = -6794

# Row 134
$ws.Range("H134").Value = 1054.7646
$ws.Range("I134").Value = 1080.1482
$ws.Range("K134").Value = 3240.4446
$ws.Range("M134").Value = -705.4446000000003

# Row 136
$ws.Range("H136").Value = 634.44446
$ws.Range("J136").Value = 691.6111
$ws.Range("L136").Value = 2074.8333
$ws.Range("N136").Value = -7174.8333

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 312.3846
$ws.Range("I33").Value = 226.83333
$ws.Range("J33").Value = 385.7143
$ws.Range("K33").Value = 1360.99998
$ws.Range("L33").Value = 2314.2858
$ws.Range("M33").Value = -1077.99998
$ws.Range("N33").Value = -2880.2858

# Row 49
$ws.Range("H49").Value = 2602
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2602
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 7806
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -8118

# Row 129
$ws.Range("H129").Value = 14369178
$ws.Range("J129").Value = 3969960.5
$ws.Range("L129").Value = 11909881.5
$ws.Range("N129").Value = -11919881.5

# Row 131
$ws.Range("H131").Value = 25001278
$ws.Range("J131").Value = 1779.1482
$ws.Range("L131").Value = 5337.444600000001
$ws.Range("N131").Value = -15417.4446

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 75002780
$ws.Range("I70").Value = 83336000
$ws.Range("K70").Value = 83336000
$ws.Range("M70").Value = -83335730

# Row 73
$ws.Range("H73").Value = 75002780
$ws.Range("I73").Value = 83336000
$ws.Range("K73").Value = 83336000
$ws.Range("M73").Value = -83335064

# Row 132
$ws.Range("H132").Value = 2273.5833
$ws.Range("I132").Value = 1657.25
$ws.Range("K132").Value = 4971.75
$ws.Range("M132").Value = -2441.75

$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 40
$ws.Range("H40").Value = 3011.2222
$ws.Range("I40").Value = 2908.2
$ws.Range("J40").Value = 3140
$ws.Range("K40").Value = 2908.2
$ws.Range("L40").Value = 3140
$ws.Range("M40").Value = -2772.2
$ws.Range("N40").Value = -3412

# Row 46
$ws.Range("H46").Value = 1800
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2333.3333
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2333.3333
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -2709.3333

# Row 93
$ws.Range("H93").Value = 611.7
$ws.Range("J93").Value = 749.5
$ws.Range("L93").Value = 749.5
$ws.Range("N93").Value = -3245.5

# Row 132
$ws.Range("H132").Value = 23274.639
$ws.Range("I132").Value = 1560.9131
$ws.Range("J132").Value = 44083.625
$ws.Range("K132").Value = 4682.7393
$ws.Range("L132").Value = 132250.875
$ws.Range("M132").Value = -2152.7393
$ws.Range("N132").Value = -137310.875

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1038.9
$ws.Range("I100").Value = 1308.3334
$ws.Range("K100").Value = 2616.6668
$ws.Range("M100").Value = -2075.6668

# Row 113
$ws.Range("H113").Value = 507.54544
$ws.Range("I113").Value = 372.875
$ws.Range("K113").Value = 1118.625
$ws.Range("M113").Value = 1051.375

# Row 126
$ws.Range("H126").Value = 111113144
$ws.Range("I126").Value = 200001760
$ws.Range("J126").Value = 2375
$ws.Range("K126").Value = 600005280
$ws.Range("L126").Value = 7125
$ws.Range("M126").Value = -600002810
$ws.Range("N126").Value = -12065
